# removed copyright to avoid politics
# - Bump the "Date" metadata value.
# - Replace the "Contact" metadata value (drop the old ContactDetail note,
#   add the MITRE/GitHub attribution instead).
# - Insert a new "Jurisdiction" row right after "Contact" (pushing
#   Description/Purpose/Copyright/Immutable down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Give the new last row (15) the same look (borders/alignment) as row 14
# before we shuffle any values into it.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)

# Shift rows 11..14 down into 12..15, bottom-up so we never clobber a row
# before reading it.
for ($r = 14; $r -ge 11; $r--) {
    $destRow = $r + 1
    $ws.Cells.Item($destRow, 1).Value = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($destRow, 2).Value = $ws.Cells.Item($r, 2).Value()
}

# Row 11 becomes the new "Jurisdiction" row (value intentionally blank).
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Update Date and Contact values.
$ws.Range("B8").Value = "2024-09-09T14:48:24-05:00"
$ws.Range("B10").Value = "MITRE, Inc (https://github.com/awatson1978/us-state-profiles)"
